$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Wrap the <fr>...</fr> run pair that precedes "desramonet" with <df>...</df>,
#    and drop the comment ("Definition!") anchored on the final "t" of that word.
# ---------------------------------------------------------------------------

# Locate the unique anchor text in the document body.
$text = $d.Content.Text
$idx = $text.IndexOf("desramonet")

if ($idx -ge 0) {
    # The 4 characters immediately before "desramonet" are the "<fr>" run.
    $rngOpen = $d.Range($idx - 4, $idx)
    if ($rngOpen.Text -eq "<fr>") {
        $rngOpen.Text = "<df><fr>"
    }

    # Re-resolve the anchor (its position does not move, but be safe) and grab
    # the 5 characters right after "desramonet", which are the "</fr>" run.
    $text = $d.Content.Text
    $idx = $text.IndexOf("desramonet")
    $rngClose = $d.Range($idx + 10, $idx + 15)
    if ($rngClose.Text -eq "</fr>") {
        $rngClose.Text = "</fr></df>"
    }
}

# Remove the "Definition!" comment left on the word "desramonet" (this also
# clears the commentRangeStart/commentRangeEnd/commentReference markers).
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}

# ---------------------------------------------------------------------------
# 2. Make the (previously implicit) footer distance explicit: 36pt = 720 twips.
# ---------------------------------------------------------------------------
$d.PageSetup.FooterDistance = 36
